# Actualización automática 2025-08-27 10:00:10
# Apply updated sales / compliance figures to the three report sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "VENTAS POR GRUPO"
# ---------------------------------------------------------------------------
$wsGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")

$wsGrupo.Range("L43").Value = 295.63

$wsGrupo.Range("L46").Value = 159.22
$wsGrupo.Range("M46").Value = 2352.19

$wsGrupo.Range("H47").Value = 566.1
$wsGrupo.Range("I47").Value = 243.9

$wsGrupo.Range("H57").Value = "2 de 55"
$wsGrupo.Range("I57").Value = "4 de 55"
$wsGrupo.Range("L57").Value = "5 de 55"

# ---------------------------------------------------------------------------
# Sheet "VENTA MENSUAL"
# ---------------------------------------------------------------------------
$wsMensual = $wb.Worksheets.Item("VENTA MENSUAL")

$wsMensual.Range("F43").Value = 295.63
$wsMensual.Range("F46").Value = 2511.41
$wsMensual.Range("F47").Value = 1549.1
$wsMensual.Range("F57").Value = 26398.26

# ---------------------------------------------------------------------------
# Sheet "CUMPLIMIENTO MENSUAL"
# ---------------------------------------------------------------------------
$wsCumpl = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

$wsCumpl.Range("D7").Value = 637.2
$wsCumpl.Range("E7").Value = 2562.8
$wsCumpl.Range("F7").Value = 0.199125

$wsCumpl.Range("D8").Value = 399.9
$wsCumpl.Range("E8").Value = 600.1
$wsCumpl.Range("F8").Value = 0.3999

$wsCumpl.Range("D15").Value = 5152.26
$wsCumpl.Range("E15").Value = 15537.74
$wsCumpl.Range("F15").Value = 0.2490217496375061

$wsCumpl.Range("D16").Value = 16038.31
$wsCumpl.Range("E16").Value = 42682.92000000001
$wsCumpl.Range("F16").Value = 0.2731262611495024

$wsCumpl.Range("D19").Value = 41564.77
$wsCumpl.Range("E19").Value = 68303.98000000001
$wsCumpl.Range("F19").Value = 0.3783129415780193

Write-Host "Actualizacion automatica aplicada."
